# Applies the dual-variable update for iteration 8 (branch and price with
# L-shaped subproblems) across the three worksheets: u_MAB, u_EOH, v_l.

$wb = $excel.ActiveWorkbook

# --- Sheet: u_MAB ---
$wsMAB = $wb.Worksheets.Item("u_MAB")

$wsMAB.Range("B15").Value = 0
$wsMAB.Range("A16").Value = 0
$wsMAB.Range("B16").Value = 0.2831078632940902
$wsMAB.Range("B24").Value = 2.183573134663976
$wsMAB.Range("A40").Value = 0
$wsMAB.Range("B40").Value = 0
$wsMAB.Range("A49").Value = 0.2264131580256175
$wsMAB.Range("B49").Value = 0.4612574414132058
$wsMAB.Range("B50").Value = 0.1942236899610495
$wsMAB.Range("B51").Value = 0.04104491635856444
$wsMAB.Range("A52").Value = 0.05182702263477301
$wsMAB.Range("B61").Value = 0

# --- Sheet: u_EOH ---
$wsEOH = $wb.Worksheets.Item("u_EOH")

$wsEOH.Range("A2").Value = -0.3457160290547003
$wsEOH.Range("A3").Value = -0.2700510980523144

# --- Sheet: v_l ---
$wsVL = $wb.Worksheets.Item("v_l")

$wsVL.Range("A2").Value = 1935936.426021874
$wsVL.Range("A3").Value = 2583294.482133841
$wsVL.Range("A4").Value = 0
